$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Overview" (sheet1): add row 3 for the new handoff file
#   b322f83c-41b1-4f58-bc25-a2d50521ef37.md
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "b322f83c-41b1-4f58-bc25-a2d50521ef37.md"
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/2e37467867a5b9008294dcce2aa64e495f00ea75/e2e/b322f83c-41b1-4f58-bc25-a2d50521ef37.md",
    "",
    "",
    "b322f83c-41b1-4f58-bc25-a2d50521ef37.md"
)
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-30-12 18:30:16"

# ----------------------------------------------------------------------
# Sheet "zh-cn" (sheet2): add row 3 for the new handoff file
# ----------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A3").Value = "b322f83c-41b1-4f58-bc25-a2d50521ef37.md"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/2e37467867a5b9008294dcce2aa64e495f00ea75/e2e/b322f83c-41b1-4f58-bc25-a2d50521ef37.md",
    "",
    "",
    "b322f83c-41b1-4f58-bc25-a2d50521ef37.md"
)

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/2e37467867a5b9008294dcce2aa64e495f00ea75/e2e/b322f83c-41b1-4f58-bc25-a2d50521ef37.md",
    "",
    "",
    ".md"
)

$wsZhCn.Range("C3").Value = "Ready for handoff"

$wsZhCn.Range("D3").Value = "b322f83c-41b1-4f58-bc25-a2d50521ef37.0c519d2e57d594fdba9fafb6f4e690b450582c1e.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/55e6a6b5e66f1f0c6cfbf332a5201b6c54e68bb4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b322f83c-41b1-4f58-bc25-a2d50521ef37.0c519d2e57d594fdba9fafb6f4e690b450582c1e.zh-cn.xlf",
    "",
    "",
    "b322f83c-41b1-4f58-bc25-a2d50521ef37.0c519d2e57d594fdba9fafb6f4e690b450582c1e.zh-cn.xlf"
)

$wsZhCn.Range("E3").Value = "2016-03-12 18:30:07"
$wsZhCn.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("H3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I3").Value = "Include"

# ----------------------------------------------------------------------
# Sheet "de-de" (sheet3): add row 3 for the new handoff file
# ----------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A3").Value = "b322f83c-41b1-4f58-bc25-a2d50521ef37.md"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/2e37467867a5b9008294dcce2aa64e495f00ea75/e2e/b322f83c-41b1-4f58-bc25-a2d50521ef37.md",
    "",
    "",
    "b322f83c-41b1-4f58-bc25-a2d50521ef37.md"
)

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/2e37467867a5b9008294dcce2aa64e495f00ea75/e2e/b322f83c-41b1-4f58-bc25-a2d50521ef37.md",
    "",
    "",
    ".md"
)

$wsDeDe.Range("C3").Value = "Ready for handoff"

$wsDeDe.Range("D3").Value = "b322f83c-41b1-4f58-bc25-a2d50521ef37.0c519d2e57d594fdba9fafb6f4e690b450582c1e.de-de.xlf"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/19a826b9f70a90bb35b1c454e1a801de3fbcdee9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b322f83c-41b1-4f58-bc25-a2d50521ef37.0c519d2e57d594fdba9fafb6f4e690b450582c1e.de-de.xlf",
    "",
    "",
    "b322f83c-41b1-4f58-bc25-a2d50521ef37.0c519d2e57d594fdba9fafb6f4e690b450582c1e.de-de.xlf"
)

$wsDeDe.Range("E3").Value = "2016-03-12 18:30:16"
$wsDeDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I3").Value = "Include"
